$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("list")
$ws1.Range("E2").Value = 1440
$ws1.Range("F2").Value = 2.93
$ws1.Range("I2").Value = 4
$ws1.Range("E3").Value = 1431
$ws1.Range("F3").Value = 3.27
$ws1.Range("I3").Value = 4.6
$ws1.Range("E4").Value = 1430
$ws1.Range("F4").Value = 2.93
$ws1.Range("H4").Value = 0.33
$ws1.Range("I4").Value = 4.67
$ws1.Range("E5").Value = 1428
$ws1.Range("F5").Value = 3.13
$ws1.Range("H5").Value = 0.6
$ws1.Range("I5").Value = 4.8
$ws1.Range("E6").Value = 1429
$ws1.Range("F6").Value = 3
$ws1.Range("G6").Value = 0.8
$ws1.Range("H6").Value = 0.93
$ws1.Range("I6").Value = 4.73
$ws1.Range("E7").Value = 1425
$ws1.Range("F7").Value = 2.67
$ws1.Range("G7").Value = 0.87
$ws1.Range("H7").Value = 1.47
$ws1.Range("I7").Value = 5
$ws1.Range("E8").Value = 1415
$ws1.Range("F8").Value = 2.33
$ws1.Range("G8").Value = 1.27
$ws1.Range("H8").Value = 2.07
$ws1.Range("I8").Value = 5.67
$ws1.Range("E9").Value = 1397
$ws1.Range("F9").Value = 3.07
$ws1.Range("G9").Value = 1.13
$ws1.Range("H9").Value = 2.67
$ws1.Range("I9").Value = 6.87
$ws1.Range("E10").Value = 1388
$ws1.Range("F10").Value = 2.67
$ws1.Range("G10").Value = 1.07
$ws1.Range("H10").Value = 3.73
$ws1.Range("I10").Value = 7.47
$ws1.Range("E11").Value = 1349
$ws1.Range("F11").Value = 2.93
$ws1.Range("G11").Value = 0.93
$ws1.Range("H11").Value = 6.2
$ws1.Range("I11").Value = 10.07
$ws1.Range("E12").Value = 1308
$ws1.Range("F12").Value = 2.4
$ws1.Range("H12").Value = 9.33
$ws1.Range("I12").Value = 12.8
$ws1.Range("E13").Value = 1246
$ws1.Range("F13").Value = 2.6
$ws1.Range("G13").Value = 1.33
$ws1.Range("H13").Value = 13
$ws1.Range("I13").Value = 16.93
$ws1.Range("E14").Value = 1178
$ws1.Range("F14").Value = 2.47
$ws1.Range("G14").Value = 0.67
$ws1.Range("H14").Value = 18.33
$ws1.Range("I14").Value = 21.47
$ws1.Range("E15").Value = 1085
$ws1.Range("F15").Value = 2.13
$ws1.Range("G15").Value = 0.6
$ws1.Range("H15").Value = 24.93
$ws1.Range("I15").Value = 27.67
$ws1.Range("E16").Value = 946
$ws1.Range("F16").Value = 2.2
$ws1.Range("G16").Value = 0.53
$ws1.Range("H16").Value = 34.2
$ws1.Range("I16").Value = 36.93
$ws1.Range("E17").Value = 763
$ws1.Range("F17").Value = 1.73
$ws1.Range("H17").Value = 47.13
$ws1.Range("I17").Value = 49.13
$ws1.Range("E18").Value = 460
$ws1.Range("F18").Value = 1
$ws1.Range("G18").Value = 0.2
$ws1.Range("H18").Value = 68.13
$ws1.Range("I18").Value = 69.33

$ws2 = $wb.Worksheets.Item("summary")
$ws2.Range("B3").Value = 1242.24
$ws2.Range("C3").Value = 280.82
$ws2.Range("D3").Value = 1388
$ws2.Range("E3").Value = 460
$ws2.Range("F3").Value = 1440
$ws2.Range("B4").Value = 2.56
$ws2.Range("C4").Value = 0.57
$ws2.Range("D4").Value = 2.67
$ws2.Range("E4").Value = 1
$ws2.Range("F4").Value = 3.27
$ws2.Range("B5").Value = 0.92
$ws2.Range("C5").Value = 0.36
$ws2.Range("E5").Value = 0.2
$ws2.Range("B6").Value = 13.71
$ws2.Range("C6").Value = 19.56
$ws2.Range("D6").Value = 3.73
$ws2.Range("F6").Value = 68.13
$ws2.Range("B7").Value = 17.18
$ws2.Range("C7").Value = 18.72
$ws2.Range("D7").Value = 7.47
$ws2.Range("E7").Value = 4
$ws2.Range("F7").Value = 69.33
